# [WIP] remove references to earlier project with tutorials
#
# The sheet had a question-less label in A2 ("Phase at Premises"); the
# author turned it into a question by appending "?". Re-assigning the
# cell's value causes the now-unused old shared string to drop out of the
# shared-strings table and the new text to be appended, which reshuffles
# the shared-string indices used by the other label cells (B2 "Single
# Phase" and A4 "power ") without altering their displayed text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Phase at Premises?"

# Author's last selection/cursor position ended up on A2 instead of A5.
$ws.Range("A2").Select()
